# Updated symbol list (GitHub Actions price refresh): column D ("Price")
# holds text values (e.g. "242.62") rather than numbers, so the raw-text
# form - including trailing zeros such as "0.001580" - must be preserved
# exactly. Writing a numeric-looking string straight into .Value lets the
# engine auto-coerce it to a number (losing formatting / trailing zeros),
# so each target cell is briefly switched to the "@" (Text) number format,
# written as a string, then restored to the "Normal" style so no lasting
# style/format change is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,23,24,25,40,41,42,43,45,47,49,50,51)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "242.88"
$ws.Cells.Item(3, 4).Value = "22.98"
$ws.Cells.Item(4, 4).Value = "5.373"
$ws.Cells.Item(5, 4).Value = "0.05965"
$ws.Cells.Item(6, 4).Value = "3.398"
$ws.Cells.Item(7, 4).Value = "6.483"
$ws.Cells.Item(8, 4).Value = "0.8061"
$ws.Cells.Item(9, 4).Value = "0.9081"
$ws.Cells.Item(10, 4).Value = "0.1418"
$ws.Cells.Item(11, 4).Value = "0.07413"
$ws.Cells.Item(12, 4).Value = "0.03311"
$ws.Cells.Item(13, 4).Value = "0.03066"
$ws.Cells.Item(14, 4).Value = "0.09343"
$ws.Cells.Item(15, 4).Value = "3.851"
$ws.Cells.Item(16, 4).Value = "0.001584"
$ws.Cells.Item(17, 4).Value = "0.04532"
$ws.Cells.Item(18, 4).Value = "0.0005935"
$ws.Cells.Item(19, 4).Value = "0.006097"
$ws.Cells.Item(20, 4).Value = "0.005020"
$ws.Cells.Item(22, 4).Value = "0.0009834"
$ws.Cells.Item(23, 4).Value = "0.00007793"
$ws.Cells.Item(24, 4).Value = "3.615"
$ws.Cells.Item(25, 4).Value = "2.139"
$ws.Cells.Item(40, 4).Value = "0.03889"
$ws.Cells.Item(41, 4).Value = "0.006077"
$ws.Cells.Item(42, 4).Value = "0.1072"
$ws.Cells.Item(43, 4).Value = "0.002692"
$ws.Cells.Item(45, 4).Value = "0.00005180"
$ws.Cells.Item(47, 4).Value = "0.0005795"
$ws.Cells.Item(49, 4).Value = "0.002259"
$ws.Cells.Item(50, 4).Value = "0.00002098"
$ws.Cells.Item(51, 4).Value = "0.0001999"

foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}

Write-Host "Updated $($rows.Count) price cells"